# Add a new entry (Purvesh Borkar / pborkar@eshopworld.com) to the
# "Distribution List" sheet, mirroring the formatting of the existing rows,
# and leave the workbook's view state (active sheet/selection) the way the
# author left it when they saved.

$wb = $excel.ActiveWorkbook

$distList = $wb.Worksheets.Item("Distribution List")
$sheet1   = $wb.Worksheets.Item("Sheet1")

# --- Append the new row to the Distribution List sheet -----------------
$distList.Range("B4").Value = "Purvesh Borkar"

# Adding the hyperlink first gives the A4 cell its link + an Excel-generated
# "Hyperlink" look; we then line the cell's style/format up with the other
# two hyperlink rows above it (A2/A3) so it matches the rest of the column,
# and make sure the visible text is the email address (not the mailto: URL).
$distList.Hyperlinks.Add($distList.Range("A4"), "mailto:pborkar@eshopworld.com") | Out-Null
$distList.Range("A3").Copy()
$distList.Range("A4").PasteSpecial(-4122)
$distList.Range("A4").Value = "pborkar@eshopworld.com"

# --- Sheet1's selection moved to A4:B4 (active cell B4) -----------------
$sheet1.Activate()
$sheet1.Range("A4:B4").Select()

# --- Distribution List is the sheet left active/selected, cursor on E7 --
$distList.Activate()
$distList.Range("E7").Select()
